# Updates cryptos list values (price + volume) to reflect latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.960.16"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.919.05"
$ws.Range("E3").Value = "  +1.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.27%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.51"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.15%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07745"
$ws.Range("E9").Value = "  +0.30%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9759"
$ws.Range("E10").Value = "  -0.53%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.47"
$ws.Range("E11").Value = "  +1.95%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.942.88"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.695"
$ws.Range("E13").Value = "  +0.43%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.957"
$ws.Range("E14").Value = "  -0.02%  "

# Row 15 - TRON
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06986"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16 - BinanceUSD (was Litecoin)
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.35%  "

# Row 17 - Litecoin (was BinanceUSD)
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.62"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009473"
$ws.Range("E18").Value = "  -0.74%  "

# Row 19 - Avalanche
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.66"
$ws.Range("E19").Value = "  -0.41%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.08%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "28.960.72"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.337"
$ws.Range("E22").Value = "  +0.19%  "

# Row 23 - Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.06"
$ws.Range("E23").Value = "  +1.58%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.179.97"
$ws.Range("E24").Value = "  -0.25%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.059"
$ws.Range("E25").Value = "  -1.54%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.94"

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.609"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.60"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.835"
$ws.Range("E30").Value = "  -0.51%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09278"
$ws.Range("E31").Value = "  +0.10%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8626"
$ws.Range("E32").Value = "  +0.07%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.105"
$ws.Range("E33").Value = "  +0.23%  "

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.238"
$ws.Range("E34").Value = "  -1.09%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.012"
$ws.Range("E35").Value = "  -0.05%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05683"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.151"
$ws.Range("E37").Value = "  +0.47%  "

# Row 38 - Frax
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.005"
$ws.Range("E38").Value = "  +0.25%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02041"
$ws.Range("E39").Value = "  +0.19%  "

# Row 40 - MXToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.085"
$ws.Range("E40").Value = "  +12.88%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.453"
$ws.Range("E41").Value = "  -0.35%  "

# Row 42 - TheSandbox
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5494"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43 - Algorand
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1755"
$ws.Range("E43").Value = "  +0.11%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.308"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 - PEPE
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002758"
$ws.Range("E45").Value = "  +16.30%  "

# Row 46 - RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.174"
$ws.Range("E46").Value = "  +3.66%  "

# Row 47 - Decentraland
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5171"
$ws.Range("E47").Value = "  -0.44%  "

# Row 48 - Cronos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06930"
$ws.Range("E48").Value = "  +1.66%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.16"
$ws.Range("E49").Value = "  -1.00%  "

# Row 50 - Quant
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.49"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51 - NEARProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.761"
$ws.Range("E51").Value = "  -0.72%  "
